# "Update material and add step 5"
# Strip the stray leading space from the 6 validation-step header labels
# (L1, N1, O1, P1, Q1, R1) now that a new step ("step 5" / BoardingPassValidation
# column) is being wired into the manifest header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L1").Value = "DateofBirth"
$ws.Range("N1").Value = "DoBValidation"
$ws.Range("O1").Value = "PersonValidation"
$ws.Range("P1").Value = "LuggageValidation"
$ws.Range("Q1").Value = "NameValidation"
$ws.Range("R1").Value = "BoardingPassValidation"

# Scroll the viewport so column F becomes the left-most visible column and
# land the selection on R2 (the new BoardingPassValidation data cell for the
# first passenger), matching where the author left off editing.
$win = $excel.ActiveWindow
$win.ScrollColumn = 6
$win.ScrollRow = 1
$ws.Range("R2").Select()
